# Code-smell bullet list: each item's text is shifted down to the next
# item (the first item's complex "Long Method (functions > 100 lines, )"
# text is removed entirely and its indent level bumped to 1; "Large
# Class (...)" becomes the simplified "Long Method"; every following
# item inherits the text that used to belong to the item above it, up
# through "Callback Hell" -> "Conditional Complexity"). "Lazy Element"
# and everything after is untouched.

$d = $word.ActiveDocument

function Get-ParagraphByText([string]$text) {
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

function Set-ParagraphTextClean([string]$oldText, [string]$newText) {
    # Locate the paragraph holding $oldText (may be spread across many
    # runs, interspersed with proofErr markers). Insert a brand-new,
    # single-run paragraph carrying the same paragraph properties right
    # before it, give that new paragraph the desired text, then remove
    # the entire original (now-shifted-to-Next) paragraph, runs,
    # proofErr markers and all.
    $target = Get-ParagraphByText $oldText
    if ($null -eq $target) {
        throw "paragraph with text [$oldText] not found"
    }
    $target.Range.InsertParagraphBefore() | Out-Null
    $oldPara = $target.Next()
    $oldPara.Range.Delete()
    if ($newText -ne "") {
        $target.Range.Text = $newText
    }
    return $target
}

# 1. "Long Method (functions > 100 lines, )" -> emptied, ilvl 0 -> 1
$p1 = Set-ParagraphTextClean "Long Method (functions > 100 lines, )" ""
$p1.Range.ListFormat.ListLevelNumber = 2

# 2. "Large Class (LoC > 200, NoA + NoM > 40)" -> "Long Method"
Set-ParagraphTextClean "Large Class (LoC > 200, NoA + NoM > 40)" "Long Method" | Out-Null

# 3. "Duplicate Code" -> "Large Class"
Set-ParagraphTextClean "Duplicate Code" "Large Class" | Out-Null

# 4. "Dead Code" -> "Duplicate Code"
Set-ParagraphTextClean "Dead Code" "Duplicate Code" | Out-Null

# 5. "Feature Envy" -> "Dead Code"
Set-ParagraphTextClean "Feature Envy" "Dead Code" | Out-Null

# 6. "Inappropriate Intimacy " -> "Feature Envy"
Set-ParagraphTextClean "Inappropriate Intimacy " "Feature Envy" | Out-Null

# 7. "Conditional Complexity" -> "Inappropriate Intimacy " (trailing space kept)
Set-ParagraphTextClean "Conditional Complexity" "Inappropriate Intimacy " | Out-Null

# 8. "Callback Hell" -> "Conditional Complexity"
Set-ParagraphTextClean "Callback Hell" "Conditional Complexity" | Out-Null

Write-Output "done"
